# Updating the proposal and CBA
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update cost figures (column C) on the Development Cost section ---
$ws.Range("C6").Value = 20000
$ws.Range("C7").Value = 10000
$ws.Range("C8").Value = 5000
$ws.Range("C9").Value = 5000
$ws.Range("C10").Value = 7500

# D9 and D10 used to be hard literal values; make them mirror D6:D8 by
# turning them into formulas that reference column C directly.
$ws.Range("D9").Formula = "= C9"
$ws.Range("D10").Formula = "= C10"

# --- Sheet view / selection changes ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C23").Select()

# --- Workbook window size / position ---
$excel.ActiveWindow.Left = 10580
$excel.ActiveWindow.Top = 710
$excel.ActiveWindow.Width = 10710
$excel.ActiveWindow.Height = 9970
